$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 3-5 with new walk-forward metric values
$ws.Range("B3").Value = 0.252
$ws.Range("C3").Value = -0.079
$ws.Range("D3").Value = 0.474
$ws.Range("E3").Value = 0.6879999999999999
$ws.Range("F3").Value = 0.786
$ws.Range("G3").Value = 0.5649999999999999

$ws.Range("B4").Value = 0.119
$ws.Range("C4").Value = -0.271
$ws.Range("D4").Value = 0.5590000000000001
$ws.Range("E4").Value = 0.748
$ws.Range("F4").Value = 0.771
$ws.Range("G4").Value = 0.467

$ws.Range("B5").Value = 0.016
$ws.Range("C5").Value = -0.216
$ws.Range("D5").Value = 0.5610000000000001
$ws.Range("E5").Value = 0.749
$ws.Range("F5").Value = 0.731
$ws.Range("G5").Value = 0.452

# Add new "Ensemble" row (row 6) with its metrics
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "Ensemble"

$ws.Range("B6").Value = 0.293
$ws.Range("C6").Value = -0.02
$ws.Range("D6").Value = 0.448
$ws.Range("E6").Value = 0.669
$ws.Range("F6").Value = 0.697
$ws.Range("G6").Value = 0.544
